$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 13.06572405778523
$ws.Range("C2").Value = 12.232762440473394
$ws.Range("D2").Value = 13.858941247051874
$ws.Range("E2").Value = 13.066501868432459

# Row 3 values
$ws.Range("B3").Value = 13.20096127938978
$ws.Range("C3").Value = 11.090405638127915
$ws.Range("D3").Value = 15.363611003681616
$ws.Range("E3").Value = 12.552884336383746

$ws.Range("B1:E3").Select()
